# Fruta / hortaliza, semanal
# Insert a new data row at row 80 (pushing existing rows 80-125 down to 81-126)
# and populate it with the latest weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 80..125 down by one to make room for the new observation.
$ws.Rows.Item(80).Insert()

# Populate the newly inserted row 80 with the new weekly record.
$ws.Range("A80").Value = 10
$ws.Range("B80").Value = "Vega Modelo de Temuco"
$ws.Range("C80").Value = "La Araucanía"
$ws.Range("D80").Value = 44438
$ws.Range("E80").Value = 9
$ws.Range("F80").Value = "Fruta"
$ws.Range("G80").Value = 100102
$ws.Range("H80").Value = "Cítricos"
$ws.Range("I80").Value = 100102006
$ws.Range("J80").Value = "Pomelo"
$ws.Range("K80").Value = "Start Ruby"
$ws.Range("L80").Value = "Primera"
$ws.Range("M80").Value = 115
$ws.Range("N80").Value = 11000
$ws.Range("O80").Value = 12000
$ws.Range("P80").Value = 11565
$ws.Range("Q80").Value = "$/caja 14 kilos granel"
$ws.Range("R80").Value = "Región de O'Higgins"
$ws.Range("S80").Value = 826
$ws.Range("T80").Value = 14
